$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix header text in B1 ("Surveys" -> "Surveyss")
$ws.Range("B1").Value = "Surveyss"

# Update the selected cell / active cell on the sheet
$ws.Range("C3").Select()

# Shrink the row height of rows 10 and 11 back to the default data-row height
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15
